{"js": "// Change \"Provide better seat availability search.\" to\n// \"Provide better seat type search.\" (commit: \"updated seat availability to\n// seat type\"). Only this one occurrence of \"seat availability\" (inside the\n// \"Provide better seat availability search.\" bullet) should change \u2014 the\n// unrelated \"...different price range, seat availability, and star rating...\"\n// sentence elsewhere in the document must stay untouched, so we search for\n// the specific, longer phrase rather than the bare word \"availability\".\n\nconst body = context.document.body;\n\nconst results = body.search(\"Provide better seat availability search.\", {\n  matchCase: true,\n  matchWholeWord: false\n});\nresults.load(\"text\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error(\"Target sentence not found: 'Provide better seat availability search.'\");\n}\n\nfor (let i = 0; i < results.items.length; i++) {\n  results.items[i].insertText(\"Provide better seat type search.\", \"Replace\");\n}\n\nawait context.sync();\n", "ps1": "# Change \"Provide better seat availability search.\" to\n# \"Provide better seat type search.\" (commit: \"updated seat availability to\n# seat type\"). Only the one occurrence of that exact sentence should change;\n# an unrelated mention of \"seat availability\" elsewhere in the document (in\n# the \"...different price range, seat availability, and star rating...\"\n# sentence) must stay untouched, so we search for the full, specific\n# sentence rather than the bare word \"availability\".\n\n$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n\n$find.Text = \"Provide better seat availability search.\"\n$find.Replacement.Text = \"Provide better seat type search.\"\n$find.Forward = $true\n$find.Wrap = 1            # wdFindContinue\n$find.Format = $false\n$find.MatchCase = $true\n$find.MatchWholeWord = $false\n$find.MatchWildcards = $false\n$find.MatchSoundsLike = $false\n$find.MatchAllWordForms = $false\n\n$found = $find.Execute(\n    $find.Text,\n    $find.MatchCase,\n    $find.MatchWholeWord,\n    $find.MatchWildcards,\n    $find.MatchSoundsLike,\n    $find.MatchAllWordForms,\n    $find.Forward,\n    $find.Wrap,\n    $find.Format,\n    $find.Replacement.Text,\n    2                      # wdReplaceAll\n)\n\nif (-not $found) {\n    throw \"Target sentence not found: 'Provide better seat availability search.'\"\n}\n"}
